$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 through 6 (old GUT46037, GUT65341, GUT65671 rows), keeping
# row 3 (GUT2366) which we then overwrite with what used to be row 4's data.
$ws.Rows("4:6").Delete()

$ws.Range("A3").Value = "even_MAG-GUT46037.fa"
$ws.Range("B3").Value = [double]"3.94020900926313e-10"
$ws.Range("C3").Value = [double]"1.050622368871636e-07"
$ws.Range("D3").Value = [double]"6.709240946679316e-05"
$ws.Range("E3").Value = [double]"0.001111838947307846"
$ws.Range("F3").Value = [double]"0.9988209631869677"
$ws.Range("G3").Value = [double]"0.9988209631869677"
$ws.Range("H3").Value = "g__UBA733"
$ws.Range("I3").Value = "g__UBA733"
